$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 164.44444
$ws.Range("I12").Value = 164.44444
$ws.Range("K12").Value = 164.44444
$ws.Range("M12").Value = 5.555560000000014

$ws.Range("H38").Value = 103.416664
$ws.Range("I38").Value = 103.416664
$ws.Range("K38").Value = 310.249992
$ws.Range("M38").Value = 61.75000799999998

$ws.Range("H88").Value = 3714.55
$ws.Range("J88").Value = 2779.1
$ws.Range("L88").Value = 2779.1
$ws.Range("N88").Value = -3591.1

$ws.Range("H91").Value = 3714.55
$ws.Range("J91").Value = 2779.1
$ws.Range("L91").Value = 2779.1
$ws.Range("N91").Value = -5587.1

$ws.Range("H132").Value = 34692.168
$ws.Range("I132").Value = 34692.168
$ws.Range("K132").Value = 104076.504
$ws.Range("M132").Value = -101546.504

$ws.Range("H141").Value = 5420.8
$ws.Range("I141").Value = 1937.4546
$ws.Range("K141").Value = 5812.3638
$ws.Range("M141").Value = -632.3638000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 205.83333
$ws.Range("I5").Value = 144.22223
$ws.Range("J5").Value = 390.66666
$ws.Range("K5").Value = 144.22223
$ws.Range("L5").Value = 390.66666
$ws.Range("M5").Value = -32.22223
$ws.Range("N5").Value = -614.66666

$ws.Range("H32").Value = 233712.95
$ws.Range("I32").Value = 254413.9
$ws.Range("K32").Value = 254413.9
$ws.Range("M32").Value = -254126.9

$ws.Range("H45").Value = 2741.889
$ws.Range("I45").Value = 1811
$ws.Range("K45").Value = 1811
$ws.Range("M45").Value = -1434

$ws.Range("H74").Value = 7027.0654
$ws.Range("I74").Value = 4547.722
$ws.Range("J74").Value = 15952.7
$ws.Range("K74").Value = 4547.722
$ws.Range("L74").Value = 15952.7
$ws.Range("M74").Value = -3673.722
$ws.Range("N74").Value = -17700.7

$ws.Range("H77").Value = 7027.0654
$ws.Range("I77").Value = 4547.722
$ws.Range("J77").Value = 15952.7
$ws.Range("K77").Value = 22738.61
$ws.Range("L77").Value = 79763.5
$ws.Range("M77").Value = -18370.61
$ws.Range("N77").Value = -88499.5

$ws.Range("H88").Value = 2380.24
$ws.Range("J88").Value = 2380.24
$ws.Range("L88").Value = 2380.24
$ws.Range("N88").Value = -3192.24

$ws.Range("H91").Value = 2380.24
$ws.Range("J91").Value = 2380.24
$ws.Range("L91").Value = 2380.24
$ws.Range("N91").Value = -5188.24

$ws.Range("H132").Value = 7150.645
$ws.Range("I132").Value = 6027.2856
$ws.Range("K132").Value = 18081.8568
$ws.Range("M132").Value = -15551.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 205.83333
$ws.Range("I4").Value = 144.22223
$ws.Range("J4").Value = 390.66666
$ws.Range("K4").Value = 144.22223
$ws.Range("L4").Value = 390.66666
$ws.Range("M4").Value = -29.22223
$ws.Range("N4").Value = -620.66666

$ws.Range("H22").Value = 258.66666
$ws.Range("J22").Value = 244.5
$ws.Range("L22").Value = 244.5
$ws.Range("N22").Value = -590.5

$ws.Range("H39").Value = 853
$ws.Range("J39").Value = 853
$ws.Range("L39").Value = 853
$ws.Range("N39").Value = -1631

$ws.Range("H80").Value = 470.94736
$ws.Range("I80").Value = 536.5833
$ws.Range("K80").Value = 536.5833
$ws.Range("M80").Value = 461.4167

$ws.Range("H83").Value = 470.94736
$ws.Range("I83").Value = 536.5833
$ws.Range("K83").Value = 2682.9165
$ws.Range("M83").Value = 2309.0835

$ws.Range("H107").Value = 1287.2727
$ws.Range("I107").Value = 1162.7894
$ws.Range("K107").Value = 1162.7894
$ws.Range("M107").Value = 757.2106000000001

$ws.Range("H134").Value = 8469.23
$ws.Range("I134").Value = 8758.333000000001
$ws.Range("K134").Value = 26274.999
$ws.Range("M134").Value = -23739.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 22666.666
$ws.Range("J18").Value = 22666.666
$ws.Range("L18").Value = 22666.666
$ws.Range("N18").Value = -23126.666

$ws.Range("H22").Value = 1404.7333
$ws.Range("I22").Value = 1006.1667
$ws.Range("K22").Value = 1006.1667
$ws.Range("M22").Value = -656.1667

$ws.Range("H99").Value = 2744.1428
$ws.Range("I99").Value = 2551.5
$ws.Range("K99").Value = 2551.5
$ws.Range("M99").Value = -1053.5

$ws.Range("H103").Value = 9499.333000000001
$ws.Range("I103").Value = 9499.333000000001
$ws.Range("K103").Value = 9499.333000000001
$ws.Range("M103").Value = -8327.333000000001

$ws.Range("H126").Value = 2744.1428
$ws.Range("I126").Value = 2551.5
$ws.Range("K126").Value = 7654.5
$ws.Range("M126").Value = -5184.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1388994.8
$ws.Range("I9").Value = 1084111.2
$ws.Range("J9").Value = 1663390
$ws.Range("K9").Value = 3252333.6
$ws.Range("L9").Value = 4990170
$ws.Range("M9").Value = -3252109.6
$ws.Range("N9").Value = -4990618

$ws.Range("H11").Value = 111111256
$ws.Range("I11").Value = 78.8
$ws.Range("J11").Value = 250000240
$ws.Range("K11").Value = 236.4
$ws.Range("L11").Value = 750000720
$ws.Range("M11").Value = -96.39999999999998
$ws.Range("N11").Value = -750001000

$ws.Range("H26").Value = 170.09091
$ws.Range("I26").Value = 168.57143
$ws.Range("J26").Value = 172.75
$ws.Range("K26").Value = 505.71429
$ws.Range("L26").Value = 518.25
$ws.Range("M26").Value = -217.71429
$ws.Range("N26").Value = -1094.25

$ws.Range("H34").Value = 2485.6
$ws.Range("I34").Value = 1713.3334
$ws.Range("J34").Value = 2678.6667
$ws.Range("K34").Value = 5140.0002
$ws.Range("L34").Value = 8036.000100000001
$ws.Range("M34").Value = -5056.0002
$ws.Range("N34").Value = -8204.000100000001

$ws.Range("H39").Value = 135986.84
$ws.Range("I39").Value = 200095
$ws.Range("K39").Value = 600285
$ws.Range("M39").Value = -599991

$ws.Range("H44").Value = 6999.8887
$ws.Range("J44").Value = 6999.8887
$ws.Range("L44").Value = 20999.6661
$ws.Range("N44").Value = -21795.6661

$ws.Range("H55").Value = 56002984
$ws.Range("I55").Value = 140000600
$ws.Range("J55").Value = 4574.5557
$ws.Range("K55").Value = 420001800
$ws.Range("L55").Value = 13723.6671
$ws.Range("M55").Value = -420001623
$ws.Range("N55").Value = -14077.6671

$ws.Range("H124").Value = 12555.75
$ws.Range("I124").Value = 10147.2
$ws.Range("J124").Value = 14276.143
$ws.Range("K124").Value = 30441.6
$ws.Range("L124").Value = 42828.429
$ws.Range("M124").Value = -25531.6
$ws.Range("N124").Value = -52648.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 46933.332
$ws.Range("J45").Value = 46933.332
$ws.Range("L45").Value = 46933.332
$ws.Range("N45").Value = -48051.332

$ws.Range("H70").Value = 15983.81
$ws.Range("J70").Value = 8216.714
$ws.Range("L70").Value = 8216.714
$ws.Range("N70").Value = -8756.714

$ws.Range("H73").Value = 15983.81
$ws.Range("J73").Value = 8216.714
$ws.Range("L73").Value = 8216.714
$ws.Range("N73").Value = -10088.714

$ws.Range("H113").Value = 1890.9667
$ws.Range("J113").Value = 1824.1111
$ws.Range("L113").Value = 1824.1111
$ws.Range("N113").Value = -6164.1111

$ws.Range("H132").Value = 12008.719
$ws.Range("I132").Value = 13752.926
$ws.Range("J132").Value = 2590
$ws.Range("K132").Value = 41258.778
$ws.Range("L132").Value = 7770
$ws.Range("M132").Value = -38728.778
$ws.Range("N132").Value = -12830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2216.6667
$ws.Range("I7").Value = 1325
$ws.Range("K7").Value = 1325
$ws.Range("M7").Value = -1213

$ws.Range("H22").Value = 2798.9048
$ws.Range("I22").Value = 2376.3845
$ws.Range("J22").Value = 3485.5
$ws.Range("K22").Value = 2376.3845
$ws.Range("L22").Value = 3485.5
$ws.Range("M22").Value = -2081.3845
$ws.Range("N22").Value = -4075.5

$ws.Range("H27").Value = 2798.9048
$ws.Range("I27").Value = 2376.3845
$ws.Range("J27").Value = 3485.5
$ws.Range("K27").Value = 2376.3845
$ws.Range("L27").Value = 3485.5
$ws.Range("M27").Value = -2269.3845
$ws.Range("N27").Value = -3699.5

$ws.Range("H40").Value = 10424.714
$ws.Range("I40").Value = 17999.666
$ws.Range("K40").Value = 17999.666
$ws.Range("M40").Value = -17863.666

$ws.Range("H55").Value = 1301.3611
$ws.Range("I55").Value = 1132.8
$ws.Range("J55").Value = 1421.762
$ws.Range("K55").Value = 1132.8
$ws.Range("L55").Value = 1421.762
$ws.Range("M55").Value = -959.8
$ws.Range("N55").Value = -1767.762

$ws.Range("H126").Value = 2216.6667
$ws.Range("I126").Value = 1325
$ws.Range("K126").Value = 3975
$ws.Range("M126").Value = -1505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 10000
$ws.Range("K34").Value = 10000
$ws.Range("M34").Value = -9797

$ws.Range("H122").Value = 42699.035
$ws.Range("I122").Value = 2400.2222
$ws.Range("J122").Value = 115236.9
$ws.Range("K122").Value = 7200.6666
$ws.Range("L122").Value = 345710.7
$ws.Range("M122").Value = -4750.6666
$ws.Range("N122").Value = -350610.7

$ws.Range("H132").Value = 2383.1667
$ws.Range("I132").Value = 1814.8182
$ws.Range("K132").Value = 5444.4546
$ws.Range("M132").Value = -2914.4546
